$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 12646.182
$ws.Range("I40").Value = 12666.667
$ws.Range("J40").Value = 12642.947
$ws.Range("K40").Value = 12666.667
$ws.Range("L40").Value = 12642.947
$ws.Range("M40").Value = -12491.667
$ws.Range("N40").Value = -12992.947
$ws.Range("H134").Value = 91103.336
$ws.Range("J134").Value = 91103.336
$ws.Range("L134").Value = 91103.336
$ws.Range("N134").Value = -101243.336
$ws.Range("H135").Value = 480.94116
$ws.Range("I135").Value = 480.94116
$ws.Range("K135").Value = 4328.47044
$ws.Range("M135").Value = -1793.47044
$ws.Range("H136").Value = 78920.664
$ws.Range("J136").Value = 78920.664
$ws.Range("L136").Value = 78920.664
$ws.Range("N136").Value = -89120.664
$ws.Range("H138").Value = 1697.7142
$ws.Range("I138").Value = 1263.5238
$ws.Range("J138").Value = 1914.8096
$ws.Range("K138").Value = 3790.5714
$ws.Range("L138").Value = 5744.4288
$ws.Range("M138").Value = 1349.4286
$ws.Range("N138").Value = -16024.4288
$ws.Range("H141").Value = 3143.2778
$ws.Range("I141").Value = 2739.9412
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 8219.8236
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -3039.8236
$ws.Range("N141").Value = -40360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3090.132
$ws.Range("I122").Value = 2607.081
$ws.Range("J122").Value = 4207.1875
$ws.Range("K122").Value = 7821.243
$ws.Range("L122").Value = 12621.5625
$ws.Range("M122").Value = -5371.243
$ws.Range("N122").Value = -17521.5625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 97746.664
$ws.Range("J132").Value = 97746.664
$ws.Range("L132").Value = 97746.664
$ws.Range("N132").Value = -107866.664
$ws.Range("H135").Value = 67856.664
$ws.Range("J135").Value = 67856.664
$ws.Range("L135").Value = 67856.664
$ws.Range("N135").Value = -77996.664
$ws.Range("H138").Value = 99937.664
$ws.Range("J138").Value = 99937.664
$ws.Range("L138").Value = 99937.664
$ws.Range("N138").Value = -110217.664
$ws.Range("H140").Value = 112659.7
$ws.Range("J140").Value = 69621.89
$ws.Range("L140").Value = 69621.89
$ws.Range("N140").Value = -79981.89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3374.9285
$ws.Range("I31").Value = 1418.909
$ws.Range("J31").Value = 10547
$ws.Range("K31").Value = 1418.909
$ws.Range("L31").Value = 10547
$ws.Range("M31").Value = -1123.909
$ws.Range("N31").Value = -11137
$ws.Range("H34").Value = 3374.9285
$ws.Range("I34").Value = 1418.909
$ws.Range("J34").Value = 10547
$ws.Range("K34").Value = 1418.909
$ws.Range("L34").Value = 10547
$ws.Range("M34").Value = -1216.909
$ws.Range("N34").Value = -10951
$ws.Range("H138").Value = 58939.668
$ws.Range("J138").Value = 61127.8
$ws.Range("L138").Value = 61127.8
$ws.Range("N138").Value = -71407.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 844.9091
$ws.Range("I117").Value = 786.75
$ws.Range("K117").Value = 2360.25
$ws.Range("M117").Value = 1081.75
$ws.Range("H121").Value = 2483.9285
$ws.Range("I121").Value = 655.6
$ws.Range("J121").Value = 3499.6667
$ws.Range("K121").Value = 1966.8
$ws.Range("L121").Value = 10499.0001
$ws.Range("M121").Value = -656.8000000000002
$ws.Range("N121").Value = -13119.0001
$ws.Range("H128").Value = 453683.34
$ws.Range("I128").Value = 453683.34
$ws.Range("K128").Value = 1361050.02
$ws.Range("M128").Value = -1356070.02
$ws.Range("H129").Value = 1715.0416
$ws.Range("I129").Value = 777.8
$ws.Range("J129").Value = 2384.5
$ws.Range("K129").Value = 2333.4
$ws.Range("L129").Value = 7153.5
$ws.Range("M129").Value = 2666.6
$ws.Range("N129").Value = -17153.5
$ws.Range("H138").Value = 3598.0312
$ws.Range("I138").Value = 2043.4706
$ws.Range("J138").Value = 5359.8667
$ws.Range("K138").Value = 6130.4118
$ws.Range("L138").Value = 16079.6001
$ws.Range("M138").Value = -990.4117999999999
$ws.Range("N138").Value = -26359.6001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13313.409
$ws.Range("I122").Value = 15229.765
$ws.Range("J122").Value = 6797.8
$ws.Range("K122").Value = 45689.295
$ws.Range("L122").Value = 20393.4
$ws.Range("M122").Value = -43239.295
$ws.Range("N122").Value = -25293.4
$ws.Range("H135").Value = 79991.664
$ws.Range("J135").Value = 79991.664
$ws.Range("L135").Value = 79991.664
$ws.Range("N135").Value = -90131.664
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 26800
$ws.Range("J42").Value = 33333.332
$ws.Range("L42").Value = 33333.332
$ws.Range("N42").Value = -34459.332
$ws.Range("H49").Value = 26800
$ws.Range("J49").Value = 33333.332
$ws.Range("L49").Value = 33333.332
$ws.Range("N49").Value = -33627.332
$ws.Range("H132").Value = 6788.7095
$ws.Range("I132").Value = 9662
$ws.Range("K132").Value = 28986
$ws.Range("M132").Value = -26456
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 24500
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 24500
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 24500
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -25928
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H82").Value = 44999
$ws.Range("J82").Value = 49998
$ws.Range("L82").Value = 49998
$ws.Range("N82").Value = -50764
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H85").Value = 44999
$ws.Range("J85").Value = 49998
$ws.Range("L85").Value = 49998
$ws.Range("N85").Value = -52650
$ws.Range("H122").Value = 3258.45
$ws.Range("I122").Value = 2978.1
$ws.Range("J122").Value = 3538.8
$ws.Range("K122").Value = 8934.299999999999
$ws.Range("L122").Value = 10616.4
$ws.Range("M122").Value = -6484.299999999999
$ws.Range("N122").Value = -15516.4
$ws.Range("H132").Value = 870717.8
$ws.Range("I132").Value = 978.44446
$ws.Range("J132").Value = 8698372
$ws.Range("K132").Value = 2935.33338
$ws.Range("L132").Value = 26095116
$ws.Range("M132").Value = -405.33338
$ws.Range("N132").Value = -26100176
$ws.Range("H140").Value = 98590.39999999999
$ws.Range("J140").Value = 98590.39999999999
$ws.Range("L140").Value = 98590.39999999999
$ws.Range("N140").Value = -108950.4
